$wb = $excel.ActiveWorkbook

$wsRef = $wb.Worksheets.Item("Municipio I")
$ws = $wb.Worksheets.Item("Municipio II")

# --- Data: Municipio II, columns A=Parioli, B=Flaminio, C=Salario ---
# Values are written column-by-column (A first, then B, then C) so the
# shared-string table is built in the same order as the source edit.
$colA = @(
    "Duke's",
    "Caffè Parnaso",
    "Ercoli 1928",
    "Lo Scoiattolo Ada",
    "Bar - Caffè delle Nazioni ai Parioli",
    "Gotha Roma",
    "Enoteca Parioli",
    "Bambu’s Parioli, Roma",
    "Palmerie Parioli",
    "Bar Villa Glori"
)

$colB = @(
    "Mostò",
    "Mediterraneo Ristorante e Giardino",
    "Metropolita",
    "Etablino - Caffè Due Fontane",
    "Frisó",
    "Cavatappi enoteca wine bar bistrot",
    "Enoteca Flaminio Roma",
    "Un Caffè con Te",
    "Ristorante ""Apoteca - Provviste Alimentari"" - Quartiere Flaminio, Roma",
    "Jacobà"
)

$colC = @(
    "Dolce caffè`n",
    "Molinari Antonio",
    "PAPY",
    "KABB",
    "Gallo Bar",
    "Dietro Le Quinte",
    "La vineria",
    "Sesto",
    "Della Manna",
    "New Age Cafè"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 1, 3).Value = $colC[$i]
}

# Match the font/colour already used throughout the workbook (theme
# accent colour) by copying the format from an existing styled cell.
$wsRef.Range("A1").Copy()
$ws.Range("A1:C10").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# C1 holds a value with a trailing line break, so wrap text on it.
$ws.Range("C1").WrapText = $true

# Widen the columns to fit the new, longer place names.
$ws.Cells.ColumnWidth = 30.77734375

# --- Selection / active sheet bookkeeping ---
$ws.Range("D1:D1048576").Select()
$ws.Activate()
